$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(1,1).Value = "UNIQUE2"
Write-Output $ws.Cells.Item(1,1).Value()
Write-Output $ws.Range("A1").Value()
